$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recursos")

# Insert a new column at T (pushes old T..AS -> U..AT).
# Excel auto-adjusts column widths, formulas ($AS$1 -> $AT$1), dimension, etc.
$ws.Columns("T").Insert()

# New column header T1: "B2B" (reuses existing shared string used by _Precios!K1)
$ws.Range("T1").Value = "B2B"

# T3 / U3: explicit (non-shared) formulas pulling B2B (col K) and Servicios (col L)
# base prices from _Precios, keyed off the new $AT$1 year cell.
$ws.Range("T3").Formula = "=_xlfn.XLOOKUP(`$AT`$1&LEFT(`$D3,6)&`$H3&`$J3,_Precios!`$S:`$S,_Precios!K:K,0)"
$ws.Range("U3").Formula = "=_xlfn.XLOOKUP(`$AT`$1&LEFT(`$D3,6)&`$H3&`$J3,_Precios!`$S:`$S,_Precios!L:L,0)"

# V3/W3/X3: totals now add the B2B column (U) instead of the old single "Servicios" T column.
$ws.Range("V3").Formula = "=`$Q3+`$U3"
$ws.Range("W3").Formula = "=`$R3+`$U3"
$ws.Range("X3").Formula = "=`$S3+`$U3"

# Re-apply the AutoFilter so its range grows to include the new last column (AS).
$ws.AutoFilterMode = $false
$ws.Range("A2:AS3").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name in sync with the filter range.
$fdb = $wb.Names.Item("Recursos!_FilterDatabase")
$fdb.RefersTo = "=Recursos!`$A`$2:`$AS`$3"

# Conditional formatting range also grows by one column (AR9991 -> AS9991).
$cf = $ws.Range("A3:AR9991").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A3:AS9991"))

# Make Recursos the active/selected sheet (was _Precios).
$ws.Activate()

Write-Output "done"
